$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("classFields")

$ws.Range("B4").Value = "redisService"
$ws.Range("D4").Value = "com.macro.mall.common.service.RedisService"
$ws.Range("B5").Value = "REDIS_KEY_ADMIN"
$ws.Range("B6").Value = "REDIS_DATABASE"
$ws.Range("D6").Value = "java.lang.String"
$ws.Range("B7").Value = "orderId"
$ws.Range("D7").Value = "java.lang.Long"
$ws.Range("B8").Value = "deliverySn"
$ws.Range("D8").Value = "java.lang.String"
$ws.Range("B16").Value = "orderOperateHistoryDao"
$ws.Range("D16").Value = "com.macro.mall.dao.OmsOrderOperateHistoryDao"
$ws.Range("B18").Value = "orderDao"
$ws.Range("D18").Value = "com.macro.mall.dao.OmsOrderDao"
$ws.Range("B20").Value = "historyList"
$ws.Range("B21").Value = "orderItemList"
$ws.Range("B28").Value = "handleMan"
$ws.Range("B29").Value = "id"
$ws.Range("D29").Value = "java.lang.Long"
$ws.Range("B30").Value = "createTime"
$ws.Range("D30").Value = "java.lang.String"
$ws.Range("B31").Value = "handleTime"
$ws.Range("B32").Value = "receiverKeyword"
$ws.Range("D32").Value = "java.lang.String"
$ws.Range("B33").Value = "status"
$ws.Range("D33").Value = "java.lang.Integer"
$ws.Range("B34").Value = "receiverProvince"
$ws.Range("B35").Value = "receiverName"
$ws.Range("B36").Value = "receiverDetailAddress"
$ws.Range("B37").Value = "status"
$ws.Range("D37").Value = "java.lang.Integer"
$ws.Range("B38").Value = "orderId"
$ws.Range("D38").Value = "java.lang.Long"
$ws.Range("B39").Value = "receiverPostCode"
$ws.Range("D39").Value = "java.lang.String"
$ws.Range("B40").Value = "receiverRegion"
$ws.Range("B41").Value = "receiverCity"
$ws.Range("B42").Value = "receiverPhone"
$ws.Range("D42").Value = "java.lang.String"
$ws.Range("B44").Value = "discountAmount"
$ws.Range("D44").Value = "java.math.BigDecimal"
$ws.Range("B45").Value = "status"
$ws.Range("D45").Value = "java.lang.Integer"
$ws.Range("B46").Value = "orderId"
$ws.Range("D46").Value = "java.lang.Long"
$ws.Range("B52").Value = "width"
$ws.Range("B53").Value = "mimeType"
$ws.Range("B54").Value = "height"
$ws.Range("B55").Value = "filename"
$ws.Range("B56").Value = "size"
$ws.Range("B62").Value = "callbackBody"
$ws.Range("B63").Value = "callbackBodyType"
$ws.Range("B64").Value = "showStatus"
$ws.Range("D64").Value = "java.lang.Integer"
$ws.Range("B65").Value = "bigPic"
$ws.Range("B66").Value = "firstLetter"
$ws.Range("B67").Value = "logo"
$ws.Range("D67").Value = "java.lang.String"
$ws.Range("B68").Value = "brandStory"
$ws.Range("D68").Value = "java.lang.String"
$ws.Range("B69").Value = "name"
$ws.Range("D69").Value = "java.lang.String"
$ws.Range("B70").Value = "sort"
$ws.Range("D70").Value = "java.lang.Integer"
$ws.Range("B71").Value = "factoryStatus"
$ws.Range("D71").Value = "java.lang.Integer"
$ws.Range("B72").Value = "handleNote"
$ws.Range("B74").Value = "handleMan"
$ws.Range("B76").Value = "receiveMan"
$ws.Range("B77").Value = "receiveNote"
$ws.Range("D77").Value = "java.lang.String"
$ws.Range("B78").Value = "status"
$ws.Range("D78").Value = "java.lang.Integer"
$ws.Range("B79").Value = "id"
$ws.Range("D79").Value = "java.lang.Long"
$ws.Range("B83").Value = "relationMapper"
$ws.Range("D83").Value = "com.macro.mall.mapper.SmsFlashPromotionProductRelationMapper"
$ws.Range("B84").Value = "relationDao"
$ws.Range("D84").Value = "com.macro.mall.dao.SmsFlashPromotionProductRelationDao"
$ws.Range("B87").Value = "ALIYUN_OSS_ACCESSKEYID"
$ws.Range("B88").Value = "ALIYUN_OSS_ENDPOINT"
$ws.Range("B89").Value = "ALIYUN_OSS_ACCESSKEYSECRET"
$ws.Range("B95").Value = "redisService"
$ws.Range("D95").Value = "com.macro.mall.common.service.RedisService"
$ws.Range("B96").Value = "resourceMapper"
$ws.Range("D96").Value = "com.macro.mall.mapper.UmsResourceMapper"
$ws.Range("B98").Value = "applicationName"
$ws.Range("D98").Value = "java.lang.String"
$ws.Range("B100").Value = "keyword"
$ws.Range("D100").Value = "java.lang.String"
$ws.Range("B101").Value = "productSn"
$ws.Range("D101").Value = "java.lang.String"
$ws.Range("B103").Value = "verifyStatus"
$ws.Range("D103").Value = "java.lang.Integer"
$ws.Range("B104").Value = "publishStatus"
$ws.Range("B105").Value = "productCategoryId"
$ws.Range("D105").Value = "java.lang.Long"
$ws.Range("B110").Value = "couponDao"
$ws.Range("D110").Value = "com.macro.mall.dao.SmsCouponDao"
$ws.Range("B111").Value = "productRelationMapper"
$ws.Range("D111").Value = "com.macro.mall.mapper.SmsCouponProductRelationMapper"
$ws.Range("B113").Value = "productRelationDao"
$ws.Range("D113").Value = "com.macro.mall.dao.SmsCouponProductRelationDao"
$ws.Range("B127").Value = "ALIYUN_OSS_DIR_PREFIX"
$ws.Range("D127").Value = "java.lang.String"
$ws.Range("B128").Value = "ossClient"
$ws.Range("D128").Value = "com.aliyun.oss.OSSClient"
$ws.Range("B130").Value = "ALIYUN_OSS_BUCKET_NAME"
$ws.Range("B131").Value = "ALIYUN_OSS_MAX_SIZE"
$ws.Range("D131").Value = "int"
$ws.Range("B132").Value = "ALIYUN_OSS_EXPIRE"
$ws.Range("D132").Value = "int"
$ws.Range("B133").Value = "ALIYUN_OSS_CALLBACK"
$ws.Range("D133").Value = "java.lang.String"
$ws.Range("B134").Value = "ALIYUN_OSS_ENDPOINT"
$ws.Range("B138").Value = "Action"
$ws.Range("B139").Value = "Resource"
$ws.Range("B140").Value = "Effect"
$ws.Range("B143").Value = "productMapper"
$ws.Range("D143").Value = "com.macro.mall.mapper.PmsProductMapper"
$ws.Range("B144").Value = "productCategoryDao"
$ws.Range("D144").Value = "com.macro.mall.dao.PmsProductCategoryDao"
$ws.Range("B145").Value = "productCategoryAttributeRelationDao"
$ws.Range("D145").Value = "com.macro.mall.dao.PmsProductCategoryAttributeRelationDao"
$ws.Range("B146").Value = "productCategoryAttributeRelationMapper"
$ws.Range("D146").Value = "com.macro.mall.mapper.PmsProductCategoryAttributeRelationMapper"
$ws.Range("B147").Value = "productCategoryMapper"
$ws.Range("D147").Value = "com.macro.mall.mapper.PmsProductCategoryMapper"
$ws.Range("B149").Value = "icon"
$ws.Range("B150").Value = "productAttributeIdList"
$ws.Range("D150").Value = "java.util.List"
$ws.Range("B151").Value = "name"
$ws.Range("D151").Value = "java.lang.String"
$ws.Range("B152").Value = "navStatus"
$ws.Range("D152").Value = "java.lang.Integer"
$ws.Range("B154").Value = "parentId"
$ws.Range("D154").Value = "java.lang.Long"
$ws.Range("B155").Value = "sort"
$ws.Range("D155").Value = "java.lang.Integer"
$ws.Range("B156").Value = "productUnit"
$ws.Range("D156").Value = "java.lang.String"
$ws.Range("B158").Value = "showStatus"
$ws.Range("B164").Value = "name"
$ws.Range("B165").Value = "url"
$ws.Range("B166").Value = "status"
$ws.Range("B167").Value = "sourceType"
$ws.Range("D167").Value = "java.lang.Integer"
$ws.Range("B169").Value = "orderType"
$ws.Range("D169").Value = "java.lang.Integer"
$ws.Range("B170").Value = "receiverKeyword"
$ws.Range("D170").Value = "java.lang.String"
$ws.Range("B171").Value = "createTime"
$ws.Range("D171").Value = "java.lang.String"
$ws.Range("B174").Value = "adminRoleRelationDao"
$ws.Range("D174").Value = "com.macro.mall.dao.UmsAdminRoleRelationDao"
$ws.Range("B175").Value = "adminRoleRelationMapper"
$ws.Range("D175").Value = "com.macro.mall.mapper.UmsAdminRoleRelationMapper"
$ws.Range("B176").Value = "loginLogMapper"
$ws.Range("D176").Value = "com.macro.mall.mapper.UmsAdminLoginLogMapper"
$ws.Range("B177").Value = "adminMapper"
$ws.Range("D177").Value = "com.macro.mall.mapper.UmsAdminMapper"
$ws.Range("B178").Value = "request"
$ws.Range("D178").Value = "javax.servlet.http.HttpServletRequest"
$ws.Range("B185").Value = "accessKeyId"
$ws.Range("B186").Value = "policy"
$ws.Range("B187").Value = "callback"
$ws.Range("B188").Value = "host"
$ws.Range("B190").Value = "username"
$ws.Range("B191").Value = "newPassword"
$ws.Range("B194").Value = "Principal"
$ws.Range("B195").Value = "Effect"
$ws.Range("B197").Value = "Resource"
$ws.Range("B200").Value = "subjectProductRelationList"
$ws.Range("B202").Value = "productAttributeValueList"
$ws.Range("B203").Value = "productLadderList"
$ws.Range("B204").Value = "prefrenceAreaProductRelationList"
$ws.Range("B205").Value = "memberPriceList"
$ws.Range("B206").Value = "skuStockList"
$ws.Range("B208").Value = "ENDPOINT"
$ws.Range("B209").Value = "ACCESS_KEY"
$ws.Range("D209").Value = "java.lang.String"
$ws.Range("B210").Value = "SECRET_KEY"
$ws.Range("B211").Value = "LOGGER"
$ws.Range("D211").Value = "org.slf4j.Logger"
$ws.Range("B212").Value = "BUCKET_NAME"
$ws.Range("B215").Value = "productVertifyRecordDao"
$ws.Range("D215").Value = "com.macro.mall.dao.PmsProductVertifyRecordDao"
$ws.Range("B216").Value = "productFullReductionMapper"
$ws.Range("D216").Value = "com.macro.mall.mapper.PmsProductFullReductionMapper"
$ws.Range("B217").Value = "skuStockDao"
$ws.Range("D217").Value = "com.macro.mall.dao.PmsSkuStockDao"
$ws.Range("B218").Value = "memberPriceMapper"
$ws.Range("D218").Value = "com.macro.mall.mapper.PmsMemberPriceMapper"
$ws.Range("B219").Value = "prefrenceAreaProductRelationMapper"
$ws.Range("D219").Value = "com.macro.mall.mapper.CmsPrefrenceAreaProductRelationMapper"
$ws.Range("B220").Value = "LOGGER"
$ws.Range("D220").Value = "org.slf4j.Logger"
$ws.Range("B221").Value = "productLadderDao"
$ws.Range("D221").Value = "com.macro.mall.dao.PmsProductLadderDao"
$ws.Range("B222").Value = "memberPriceDao"
$ws.Range("D222").Value = "com.macro.mall.dao.PmsMemberPriceDao"
$ws.Range("B223").Value = "prefrenceAreaProductRelationDao"
$ws.Range("D223").Value = "com.macro.mall.dao.CmsPrefrenceAreaProductRelationDao"
$ws.Range("B224").Value = "skuStockMapper"
$ws.Range("D224").Value = "com.macro.mall.mapper.PmsSkuStockMapper"
$ws.Range("B225").Value = "productAttributeValueDao"
$ws.Range("D225").Value = "com.macro.mall.dao.PmsProductAttributeValueDao"
$ws.Range("B226").Value = "subjectProductRelationDao"
$ws.Range("D226").Value = "com.macro.mall.dao.CmsSubjectProductRelationDao"
$ws.Range("B227").Value = "productDao"
$ws.Range("D227").Value = "com.macro.mall.dao.PmsProductDao"
$ws.Range("B228").Value = "productLadderMapper"
$ws.Range("D228").Value = "com.macro.mall.mapper.PmsProductLadderMapper"
$ws.Range("B229").Value = "productMapper"
$ws.Range("D229").Value = "com.macro.mall.mapper.PmsProductMapper"
$ws.Range("B230").Value = "productAttributeValueMapper"
$ws.Range("D230").Value = "com.macro.mall.mapper.PmsProductAttributeValueMapper"
$ws.Range("B231").Value = "subjectProductRelationMapper"
$ws.Range("D231").Value = "com.macro.mall.mapper.CmsSubjectProductRelationMapper"
$ws.Range("B232").Value = "productFullReductionDao"
$ws.Range("D232").Value = "com.macro.mall.dao.PmsProductFullReductionDao"
$ws.Range("B233").Value = "resourceService"
$ws.Range("D233").Value = "com.macro.mall.service.UmsResourceService"
$ws.Range("B234").Value = "roleMapper"
$ws.Range("D234").Value = "com.macro.mall.mapper.UmsRoleMapper"
$ws.Range("B235").Value = "roleResourceRelationMapper"
$ws.Range("D235").Value = "com.macro.mall.mapper.UmsRoleResourceRelationMapper"
$ws.Range("B236").Value = "roleDao"
$ws.Range("D236").Value = "com.macro.mall.dao.UmsRoleDao"
$ws.Range("B237").Value = "roleMenuRelationMapper"
$ws.Range("D237").Value = "com.macro.mall.mapper.UmsRoleMenuRelationMapper"
$ws.Range("B244").Value = "handAddStatus"
$ws.Range("B245").Value = "name"
$ws.Range("B246").Value = "filterType"
$ws.Range("B247").Value = "productAttributeCategoryId"
$ws.Range("D247").Value = "java.lang.Long"
$ws.Range("B248").Value = "type"
$ws.Range("B249").Value = "sort"
$ws.Range("B250").Value = "inputList"
$ws.Range("D250").Value = "java.lang.String"
$ws.Range("B251").Value = "searchType"
$ws.Range("B252").Value = "selectType"
$ws.Range("B253").Value = "inputType"
$ws.Range("B254").Value = "relatedStatus"
$ws.Range("B259").Value = "note"
$ws.Range("B260").Value = "password"
$ws.Range("B261").Value = "email"
$ws.Range("B262").Value = "nickName"
$ws.Range("B263").Value = "productAttributeCategoryMapper"
$ws.Range("D263").Value = "com.macro.mall.mapper.PmsProductAttributeCategoryMapper"
$ws.Range("B264").Value = "productAttributeDao"
$ws.Range("D264").Value = "com.macro.mall.dao.PmsProductAttributeDao"
